$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Accept column (F) for the two SSH/HTTP rows that used to point to
# "pastorious.rivetweb.org" so that they now point to the new target IP.
$ws.Range("F3").Value = "164.92.216.90"
$ws.Range("F4").Value = "164.92.216.90"

# Reflect the new active cell/selection used while making the edit.
$ws.Range("F4").Select() | Out-Null
